# Enemies.xlsx — add "Engineer" enemy row, retune Bulleter/Clarissa values.
#
# Target layout (Sheet1, table A2:G7):
#   Row2  Name       Health Attack Defense Speed EXP Value   (header, unchanged)
#   Row3  Bulleter   30     4      0       2     5   0.05    (Value 1 -> 0.05)
#   Row4  Engineer   100    5      0       8     20  1       (new row, inserted)
#   Row5  Drunk Guy  500    10     0       4     50  6.15    (shifted down from row4)
#   Row6  Clarissa   700    6      0       15    400 12.3    (shifted from row5, Value 20 -> 12.3)
#   Row7  Twin 1     4500   40     0       15    -   0       (shifted down from row6)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new blank row above the current row 4 (Drunk Guy), shifting
# Drunk Guy / Clarissa / Twin 1 down by one. This also pushes the sheet's
# used range / dimension out to row 7.
$ws.Rows.Item(4).Insert()

# --- Row 3: Bulleter -- only the Value column changes
$ws.Range("G3").Value = 0.05

# --- Row 4: new "Engineer" enemy
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "Engineer"
$ws.Range("B4").Value = 100
$ws.Range("C4").Value = 5
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 8
$ws.Range("F4").Value = 20
$ws.Range("G4").Value = 1

# --- Row 5: Drunk Guy (re-affirm values after the shift)
$ws.Range("A5").Value = "Drunk Guy"
$ws.Range("B5").Value = 500
$ws.Range("C5").Value = 10
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 4
$ws.Range("F5").Value = 50
$ws.Range("G5").Value = 6.15

# --- Row 6: Clarissa (Value 20 -> 12.3)
$ws.Range("A6").Value = "Clarissa"
$ws.Range("B6").Value = 700
$ws.Range("C6").Value = 6
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 15
$ws.Range("F6").Value = 400
$ws.Range("G6").Value = 12.3

# --- Row 7: Twin 1 (re-affirm values after the shift; no F value)
$ws.Range("A7").Value = "Twin 1"
$ws.Range("B7").Value = 4500
$ws.Range("C7").Value = 40
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 15
$ws.Range("G7").Value = 0

# Grow the table (ListObject) to cover the new row
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A2:G7")) | Out-Null

# Selection now sits on the last populated cell, G7
$ws.Range("G7").Select() | Out-Null
